# petty-cashBook-2021.xlsx -- "Update 26-Mei-2021, end of day update."
# Fills in the transaction rows for 25-Mei-2021 (row 3, continuing) and
# 26-Mei-2021 (from row 11) on Sheet1 of the petty cash book.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 25-Mei-2021 continued (row 3 already had the date/description) ---
$ws.Range("D3").Formula = "=45000+210000"

$ws.Range("B4").Value = "BELI kresek"
$ws.Range("D4").Formula = "=92500"

$ws.Range("B5").Value = "TRANSFER BCA"
$ws.Range("D5").Formula = "=2550000+430000+3650000+60000+419000+120000+3005000+1879000+11000000"

$ws.Range("B6").Value = "A/R"
$ws.Range("C6").Formula = "=3650000+9071000+11000000+45357000"

$ws.Range("B7").Value = "FREIGHT OUT"
$ws.Range("D7").Formula = "=218500"

$ws.Range("B8").Value = "SALES - cash/retail"
$ws.Range("C8").Formula = "=64843525-12376025-45357000"

$ws.Range("B9").Value = "SELISIH - lebih"
$ws.Range("C9").Value = 110000

$ws.Range("B10").Value = "SETOR KE BANK"
$ws.Range("D10").Formula = "=53000000"

# --- 26-Mei-2021 ---
$ws.Range("A11").Value = 44341
$ws.Range("B11").Value = "Wages Expense"
$ws.Range("D11").Formula = "=60000+240000"

$ws.Range("B12").Value = "BELI nota"
$ws.Range("D12").Formula = "=100000"

$ws.Range("B13").Value = "TRANSFER BCA"
$ws.Range("D13").Formula = "=380000+1356500"

$ws.Range("B14").Value = "JASON - visa Kanada"
$ws.Range("D14").Value = 2950000

$ws.Range("B15").Value = "A/R"
$ws.Range("C15").Formula = "=2755500"

$ws.Range("B16").Value = "SALES - cash/retail"
$ws.Range("C16").Formula = "=4102975+7102525-2755500"

$ws.Range("B17").Value = "SELISIH - kurang"
$ws.Range("D17").Value = 100000

$ws.Range("B18").Value = "SETOR KE BANK"
$ws.Range("D18").Value = 6000000

$ws.Range("A19").Value = 44342
$ws.Range("B19").Value = "Wages Expense"
$ws.Range("D19").Value = 60000

$ws.Range("B20").Value = "A/R"
$ws.Range("C20").Formula = "=40965000+8890000+21622500"

$ws.Range("B21").Value = "TRANSFER BCA"
$ws.Range("D21").Formula = "=40965000+21082500+9432000"

# --- extend the running-balance formula one more row, to the new row 114 ---
$ws.Range("E114").Formula = "=E113+C114-D114"

# leave the selection where the day's last entry was made
$ws.Range("B22").Select()
